# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.655.07"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$ws.Range("D3").Value = "1.963.70"
$ws.Range("E3").Value = "  +2.38%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.14"
$ws.Range("E5").Value = "  +1.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4840"
$ws.Range("E7").Value = "  +0.90%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2949"
$ws.Range("E8").Value = "  +1.77%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06793"
$ws.Range("E9").Value = "  +1.14%  "

# Row 10
$ws.Range("B10").Value = "Litecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.18"
$ws.Range("E10").Value = "  -0.25%  "

# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.41"
$ws.Range("E11").Value = "  +1.74%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.975.92"
$ws.Range("E12").Value = "  +3.05%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07769"
$ws.Range("E13").Value = "  +2.63%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.471"
$ws.Range("E14").Value = "  +3.82%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6893"
$ws.Range("E15").Value = "  +3.16%  "

# Row 16
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "294.90"
$ws.Range("E16").Value = "  -1.40%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.677.02"
$ws.Range("E17").Value = "  +0.53%  "

# Row 18
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.245.67"
$ws.Range("E18").Value = "  +3.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +1.86%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007716"
$ws.Range("E20").Value = "  +1.93%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.620"
$ws.Range("E21").Value = "  +0.89%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.617"
$ws.Range("E24").Value = "  +2.44%  "

# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.921"
$ws.Range("E25").Value = "  +4.59%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.56"
$ws.Range("E26").Value = "  +3.67%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.21"
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.191"
$ws.Range("E28").Value = "  +3.70%  "

# Row 29
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1066"
$ws.Range("E29").Value = "  -1.08%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.436"
$ws.Range("E30").Value = "  +2.69%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.737"
$ws.Range("E31").Value = "  +17.16%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.451"
$ws.Range("E32").Value = "  +6.83%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05118"
$ws.Range("E33").Value = "  +2.42%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7727"
$ws.Range("E34").Value = "  +4.92%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.182"
$ws.Range("E35").Value = "  +3.94%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02047"
$ws.Range("E36").Value = "  +0.58%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.732"
$ws.Range("E37").Value = "  +0.40%  "

# Row 38
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.719"
$ws.Range("E38").Value = "  +1.35%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.121"
$ws.Range("E39").Value = "  +5.03%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.362"
$ws.Range("E40").Value = "  +7.89%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4482"
$ws.Range("E41").Value = "  +1.09%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.33"
$ws.Range("E42").Value = "  -1.45%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8773"
$ws.Range("E43").Value = "  +1.79%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.43"
$ws.Range("E44").Value = "  -2.25%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.539"
$ws.Range("E46").Value = "  +3.65%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1284"
$ws.Range("E47").Value = "  +4.30%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.386"
$ws.Range("E48").Value = "  +1.55%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.14"
$ws.Range("E49").Value = "  +3.10%  "

# Row 50
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.71"
$ws.Range("E50").Value = "  -3.34%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "920.35"
$ws.Range("E51").Value = "  +6.66%  "
